$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Row 11: new changelog entry
$ws.Range("A11").Value = 41676
$ws.Range("B11").Value = "Adicionado indice na coluna email na tabela users (agora é campo único)"
$ws.Range("C11").Value = "não"

# Row 12: new changelog entry (copy date formatting from the row above)
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = 41676
$ws.Range("B12").Value = "userType na table users agora é varchar"
$ws.Range("C12").Value = "não"

$ws.Range("C12").Select()
